$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "agua-corriente" mapping file reference from F5 to B5 (keeps cell
# formatting intact, unlike writing straight into a previously-blank cell).
$ws.Range("F5").Cut($ws.Range("B5"))
$ws.Range("F5").Clear()

# Row 1 - column headers (human-readable labels, now capitalised / accented)
$ws.Range("A1").Value = "Número de edificios"
$ws.Range("B1").Value = "Agua corriente"
$ws.Range("C1").Value = "Comarca nombre"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "Provincia código"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Provincia nombre"

# Row 2 - concept identifiers (measure / dimension URIs)
$ws.Range("A2").Value = "iaest-measure:numero-de-edificios"
$ws.Range("B2").Value = "iaest-dimension:agua-corriente"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# Row 3 - role (measure vs dimension)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "dim"

# Row 4 - data type / code-list reference
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "skos:Concept"
$ws.Range("C4").Value = "URI-comarca"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "URI-Provincia"
